$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: change politeness_score (B2) from text "3" to a real number 3.
$ws.Range("B2").Value = 3

# Row 3: new annotation row for Ying Tang.
$ws.Range("A3").Value = "Ying Tang"

# B3 should stay a text "3" (like the original B2 before this edit), not a number.
$b3 = $ws.Range("B3")
$b3.NumberFormat = "@"
$b3.Value = "3"
$b3.Style = "Normal"

$ws.Range("C3").Value = "should be mentioned"
$ws.Range("D3").Value = "SUG"
$ws.Range("E3").Value = "MET"
$ws.Range("F3").Value = "c8048836-24fe-4e27-95aa-c7cfb58ac155"
$ws.Range("G3").Value = "rkc_hGb0Z_annotated.xlsx"
$ws.Range("H3").Value = "The structure of the global policies used in the experiments should be mentioned somewhere."
